$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medical Devices")

# Insert 6 new rows before row 20 to make room for new companies
$ws.Range("A20:A25").EntireRow.Insert()

# Fill the 6 new rows (A = "x" flag, B = company name, C = ticker)
$ws.Cells.Item(20, 1).Value = "x"
$ws.Cells.Item(20, 2).Value = "Idexx Labs"
$ws.Cells.Item(20, 3).Value = "IDXX"

$ws.Cells.Item(21, 1).Value = "x"
$ws.Cells.Item(21, 2).Value = "Dexcom"
$ws.Cells.Item(21, 3).Value = "DXCM"

$ws.Cells.Item(22, 1).Value = "x"
$ws.Cells.Item(22, 2).Value = "Sartorius"
$ws.Cells.Item(22, 3).Value = "SRT GR"

$ws.Cells.Item(23, 1).Value = "x"
$ws.Cells.Item(23, 2).Value = "Fujifilm"
$ws.Cells.Item(23, 3).Value = "4901 JP"

$ws.Cells.Item(24, 1).Value = "x"
$ws.Cells.Item(24, 2).Value = "Mettler-Toledo"
$ws.Cells.Item(24, 3).Value = "MTD"

$ws.Cells.Item(25, 1).Value = "x"
$ws.Cells.Item(25, 2).Value = "Olympus"
$ws.Cells.Item(25, 3).Value = "7733 JP"

# The old rows shifted down by 6: old Dexcom row (was 20) is now 26,
# old Olympus row (was 21) is now 27, and old Sartorius row (was 34) is
# now 40. Their data now lives in the new block above, so remove the
# stale duplicate rows (delete bottom-most first to keep indices stable).
$ws.Range("B40").EntireRow.Delete()
$ws.Range("B27").EntireRow.Delete()
$ws.Range("B26").EntireRow.Delete()

# Update the active-cell selection to match the saved view state
$ws.Range("B25").Select() | Out-Null
